# Refresh cryptocurrency price/volume data (and two pairs of re-ranked rows)
# as captured by the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to Text format before assigning so Excel does not
    # auto-coerce numeric-looking strings (e.g. "6.67", "0.401") into
    # real numbers, then drop the temporary format so the cell keeps
    # its original (default) style - only the stored string changes.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-TextValue "D2" '65.845.36'
Set-TextValue "E2" '  -1.42%  '
Set-TextValue "D3" '3.294.36'
Set-TextValue "E3" '  -0.25%  '
Set-TextValue "E4" '  +0.02%  '
Set-TextValue "D5" '573.43'
Set-TextValue "E5" '  -0.29%  '
Set-TextValue "D6" '178.39'
Set-TextValue "E6" '  -3.98%  '
Set-TextValue "E7" '  +4.79%  '
Set-TextValue "E8" '  +0.02%  '
Set-TextValue "E9" '  -2.65%  '
Set-TextValue "D10" '6.67'
Set-TextValue "E10" '  +0.10%  '
Set-TextValue "D11" '0.401'
Set-TextValue "E11" '  -2.34%  '
Set-TextValue "D12" '3.867.21'
Set-TextValue "E12" '  -0.20%  '
Set-TextValue "E13" '  -3.60%  '
Set-TextValue "D14" '26.58'
Set-TextValue "E14" '  -3.12%  '
Set-TextValue "D15" '65.975.95'
Set-TextValue "E15" '  -1.63%  '
Set-TextValue "B16" 'WrappedEther'
Set-TextValue "C16" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D16" '3.320.75'
Set-TextValue "E16" '  +1.02%  '
Set-TextValue "B17" 'ShibaInu'
Set-TextValue "C17" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D17" '0.0000163'
Set-TextValue "E17" '  -1.78%  '
Set-TextValue "D18" '435.07'
Set-TextValue "E18" '  -1.59%  '
Set-TextValue "E19" '  -1.92%  '
Set-TextValue "D20" '13.31'
Set-TextValue "E20" '  -1.47%  '
Set-TextValue "E21" '  -4.27%  '
Set-TextValue "D22" '72.43'
Set-TextValue "E22" '  -2.43%  '
Set-TextValue "E23" '  +0.21%  '
Set-TextValue "B24" 'Polygon'
Set-TextValue "C24" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D24" '0.514'
Set-TextValue "E24" '  +0.08%  '
Set-TextValue "B25" 'WrappedeETH'
Set-TextValue "C25" 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue "D25" '3.448.43'
Set-TextValue "E25" '  +0.40%  '
Set-TextValue "E26" '  -3.94%  '
Set-TextValue "D27" '0.195'
Set-TextValue "E27" '  +4.51%  '
Set-TextValue "D28" '8.96'
Set-TextValue "E28" '  -1.48%  '
Set-TextValue "E29" '  -0.11%  '
Set-TextValue "E30" '  -1.46%  '
Set-TextValue "D31" '22.41'
Set-TextValue "E31" '  -1.79%  '
Set-TextValue "E32" '  +0.02%  '
Set-TextValue "D33" '5.14'
Set-TextValue "E33" '  -3.49%  '
Set-TextValue "D34" '6.64'
Set-TextValue "E34" '  -1.90%  '
Set-TextValue "E35" '  -3.37%  '
Set-TextValue "E36" '  -4.78%  '
Set-TextValue "D37" '157.14'
Set-TextValue "E37" '  -3.47%  '
Set-TextValue "D38" '27.09'
Set-TextValue "E38" '  -1.59%  '
Set-TextValue "E39" '  -2.97%  '
Set-TextValue "D40" '2.781.73'
Set-TextValue "E40" '  +1.91%  '
Set-TextValue "E41" '  -0.06%  '
Set-TextValue "E42" '  -2.65%  '
Set-TextValue "D43" '40.41'
Set-TextValue "E43" '  +0.61%  '
Set-TextValue "D44" '6.11'
Set-TextValue "E44" '  -2.34%  '
Set-TextValue "E45" '  -1.85%  '
Set-TextValue "D46" '323.37'
Set-TextValue "E46" '  -1.27%  '
Set-TextValue "E47" '  -3.81%  '
Set-TextValue "E48" '  -4.71%  '
Set-TextValue "E49" '  -1.34%  '
Set-TextValue "E50" '  +2.51%  '
Set-TextValue "D51" '0.999'
Set-TextValue "E51" '  +0.07%  '
